# Applies the commit "Added yopmail component in common functions"
# Change: cell J3 on Sheet1 (Keywords2 for the testSendBusinessInvitations row)
# gets additional -p parameters appended to its BusinessProfileTest keyword string.
# The selected cell in the sheet view also moves from J3 to I3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValue = "coyni.admin.paymentstests.BusinessProfileTest,`ntestAddMerchantUser,`n-pfirstName,`n-plastName,`n-pheading,`n-pinviteHeading,`n-pCreateAccountHeading,`n-pphoneHeading,`n-pcode,`n-pemailHeading"

$ws.Range("J3").Value = $newValue

# Update the active selection to I3 (matches the saved sheet view state in the diff)
$ws.Activate()
$ws.Range("I3").Select()
